$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8823
$ws1.Range("F3").Value = 8189
$ws1.Range("F5").Value = 198
$ws1.Range("F6").Value = 39
$ws1.Range("F9").Value = 148
$ws1.Range("F11").Value = 249
$ws1.Range("F12").Value = 739
$ws1.Range("F13").Value = 201
$ws1.Range("F14").Value = 5073
$ws1.Range("F19").Value = 153
$ws1.Range("F20").Value = 137
$ws1.Range("F21").Value = 2

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8823
$ws4.Range("F3").Value = 8189
$ws4.Range("F5").Value = 198
$ws4.Range("F6").Value = 39
$ws4.Range("F9").Value = 148
$ws4.Range("F11").Value = 249
$ws4.Range("F12").Value = 739
$ws4.Range("F13").Value = 201
$ws4.Range("F14").Value = 5074
$ws4.Range("F19").Value = 153
$ws4.Range("F20").Value = 137
$ws4.Range("F21").Value = 2
